# Helper: write a value as TEXT (shared string), regardless of whether it
# looks numeric, by routing it through a formula that evaluates to a text
# result and then freezing that result as a static value via copy/paste-
# special-values. This avoids Excel's automatic "looks like a number ->
# store as number" coercion that a plain .Value assignment triggers.
function Set-TextValue {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet1: benchmark summary (A: metric name, B: score text)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Columns.Item(1).ColumnWidth = 20.7109375

$sheet1Data = @(
    @("Overall Score", "1726$([char]10)"),
    @("Productivity", "1614$([char]10)"),
    @("Creativity", "1914$([char]10)"),
    @("Responsiveness", "1549$([char]10)")
)

for ($i = 0; $i -lt $sheet1Data.Length; $i++) {
    $r = $i + 1
    Set-TextValue $ws1.Cells.Item($r, 1) $sheet1Data[$i][0]
    Set-TextValue $ws1.Cells.Item($r, 2) $sheet1Data[$i][1]
}

# ---------------------------------------------------------------
# Sheet2: benchmark category breakdown, inserted right after Sheet1
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Columns.Item(1).ColumnWidth = 20.7109375

$sheet2Data = @(
    @("Overall Score", "6223"),
    @("Essentials", "11101"),
    @("Productivity", "11101"),
    @("Digital Content Creation", "10373"),
    @("App Startup", "10373"),
    @("Video Conferencing", "5681"),
    @("Web Browsing", "5681"),
    @("Spreadsheets", "17516"),
    @("Writing", "8168"),
    @("Photo Editing", "9564"),
    @("Render and Visual", "11977"),
    @("Video Editing", "8984")
)

for ($i = 0; $i -lt $sheet2Data.Length; $i++) {
    $r = $i + 1
    Set-TextValue $ws2.Cells.Item($r, 1) $sheet2Data[$i][0]
    Set-TextValue $ws2.Cells.Item($r, 2) $sheet2Data[$i][1]
}

# Restore Sheet1 as the active/selected sheet (matches original tabSelected state)
$ws1.Activate() | Out-Null
$ws1.Cells.Item(1, 1).Select() | Out-Null
